$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coherencia")
$ws.Activate()
Write-Output "activated"
